# From v1.2 to v1.2.1
# Swap the "second step" content between TC2 and TC4 blocks in the
# "Test Suite" sheet: TC2's step becomes what used to be TC4's step, and
# TC4's step becomes what used to be TC2's step. TC3's step is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Original (before edit) text values, captured for clarity:
$tc2_action = "Chefe Dado um registro selecionado (solicitação aguardando autorização de pagamento - AP), o usuário pode atribuir/desatribuir a responsabilidade da AP a si próprio; e Clica para atribuir/desatribuir o registro a si mesmo."
$tc2_result = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela AP) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

$tc4_action = "Chefe Seleciona um usuário para filtrar as autorizações de pagamento associadas a ele; e Submete a busca ao sistema."
$tc4_result = "SYSTEM Filtra os registros (autorizações de pagamento pendentes) e exibe apenas aqueles atribuídos ao usuário selecionado."

# Row 20 holds TC2's second step (columns B = Steps, D = Expected Results).
# It should now show what used to be TC4's second step.
$ws.Range("B20").Value = $tc4_action
$ws.Range("D20").Value = $tc4_result

# Row 36 holds TC4's second step (columns B = Steps, D = Expected Results).
# It should now show what used to be TC2's second step.
$ws.Range("B36").Value = $tc2_action
$ws.Range("D36").Value = $tc2_result

# Row 28 (TC3's second step) is left unchanged.
